$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.78"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.06"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.324"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06234"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.650"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.635"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.400"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01380"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1599"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08410"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03516"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03216"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.061"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09289"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001639"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04743"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006341"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005705"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.721"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.325"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3354"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002708"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04731"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007092"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004507"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1167"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01216"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006142"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009914"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7833"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002419"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002404"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01242"
